$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 43.96295866666666
$ws.Range("H2").Value = 131.888876
$ws.Range("I2").Value = 0.1193823486802574
$ws.Range("J2").Value = 0.1297146081693155
$ws.Range("M2").Value = 0.9105153333333335
$ws.Range("N2").Value = 2.731546
$ws.Range("O2").Value = 0.03002598096739917
$ws.Range("P2").Value = 0.03049089185478593
$ws.Range("Q2").Value = 40.02894796469955
$ws.Range("R2").Value = 360.260531682296
$ws.Range("S2").Value = 0.003584572129316819
$ws.Range("T2").Value = 0.003955114089676528
$ws.Range("G3").Value = 43.96295866666666
$ws.Range("H3").Value = 131.888876
$ws.Range("I3").Value = 0.1193823486802574
$ws.Range("J3").Value = 0.1297146081693155
$ws.Range("O3").Value = 0.2937114969858886
$ws.Range("P3").Value = 0.2982592142727163
$ws.Range("Q3").Value = 391.559637709999
$ws.Range("R3").Value = 3524.036739389991
$ws.Range("S3").Value = 0.03506396834456971
$ws.Range("T3").Value = 0.0386885771122733
$ws.Range("G4").Value = 43.96295866666666
$ws.Range("H4").Value = 131.888876
$ws.Range("I4").Value = 0.1193823486802574
$ws.Range("J4").Value = 0.1297146081693155
$ws.Range("M4").Value = 9.030046333333333
$ws.Range("N4").Value = 27.090139
$ws.Range("O4").Value = 0.2977830129963756
$ws.Range("P4").Value = 0.3023937720910131
$ws.Range("Q4").Value = 396.9875537104181
$ws.Range("R4").Value = 3572.887983393764
$ws.Range("S4").Value = 0.03555003548859093
$ws.Range("T4").Value = 0.03922488965962705
$ws.Range("G5").Value = 43.96295866666666
$ws.Range("H5").Value = 131.888876
$ws.Range("I5").Value = 0.1193823486802574
$ws.Range("J5").Value = 0.1297146081693155
$ws.Range("M5").Value = 1.38711
$ws.Range("N5").Value = 2.77422
$ws.Range("O5").Value = 0.04574259975086167
$ws.Range("P5").Value = 0.03096724053022875
$ws.Range("Q5").Value = 60.98145959611999
$ws.Range("R5").Value = 365.88875757672
$ws.Range("S5").Value = 0.005460858992998822
$ws.Range("T5").Value = 0.004016903471463567
$ws.Range("G6").Value = 43.96295866666666
$ws.Range("H6").Value = 131.888876
$ws.Range("I6").Value = 0.1193823486802574
$ws.Range("J6").Value = 0.1297146081693155
$ws.Range("M6").Value = 10.089997
$ws.Range("N6").Value = 30.269991
$ws.Range("O6").Value = 0.3327369092994751
$ws.Range("P6").Value = 0.3378888812512559
$ws.Range("Q6").Value = 443.5861210577905
$ws.Range("R6").Value = 3992.275089520115
$ws.Range("S6").Value = 0.03972291372478111
$ws.Range("T6").Value = 0.04382912383627503
$ws.Range("I7").Value = 0.1955776293519722
$ws.Range("J7").Value = 0.212504409894141
$ws.Range("M7").Value = 0.9105153333333335
$ws.Range("N7").Value = 2.731546
$ws.Range("O7").Value = 0.03002598096739917
$ws.Range("P7").Value = 0.03049089185478593
$ws.Range("Q7").Value = 65.57725522185223
$ws.Range("R7").Value = 590.1952969966701
$ws.Range("S7").Value = 0.005872410176571367
$ws.Range("T7").Value = 0.006479448980747354
$ws.Range("I8").Value = 0.1955776293519722
$ws.Range("J8").Value = 0.212504409894141
$ws.Range("O8").Value = 0.2937114969858886
$ws.Range("P8").Value = 0.2982592142727163
$ws.Range("S8").Value = 0.05744339829391901
$ws.Range("T8").Value = 0.06338139832451375
$ws.Range("I9").Value = 0.1955776293519722
$ws.Range("J9").Value = 0.212504409894141
$ws.Range("M9").Value = 9.030046333333333
$ws.Range("N9").Value = 27.090139
$ws.Range("O9").Value = 0.2977830129963756
$ws.Range("P9").Value = 0.3023937720910131
$ws.Range("Q9").Value = 650.3631859754339
$ws.Range("R9").Value = 5853.268673778905
$ws.Range("S9").Value = 0.05823969574311867
$ws.Range("T9").Value = 0.06426001009386412
$ws.Range("I10").Value = 0.1955776293519722
$ws.Range("J10").Value = 0.212504409894141
$ws.Range("M10").Value = 1.38711
$ws.Range("N10").Value = 2.77422
$ws.Range("O10").Value = 0.04574259975086167
$ws.Range("P10").Value = 0.03096724053022875
$ws.Range("Q10").Value = 99.90261905615
$ws.Range("R10").Value = 599.4157143369
$ws.Range("S10").Value = 0.008946229219669638
$ws.Range("T10").Value = 0.006580675174926187
$ws.Range("I11").Value = 0.1955776293519722
$ws.Range("J11").Value = 0.212504409894141
$ws.Range("M11").Value = 10.089997
$ws.Range("N11").Value = 30.269991
$ws.Range("O11").Value = 0.3327369092994751
$ws.Range("P11").Value = 0.3378888812512559
$ws.Range("Q11").Value = 726.7030924502716
$ws.Range("R11").Value = 6540.327832052444
$ws.Range("S11").Value = 0.06507589591869353
$ws.Range("T11").Value = 0.07180287732008964
$ws.Range("G12").Value = 75.27587666666666
$ws.Range("H12").Value = 225.82763
$ws.Range("I12").Value = 0.2044132430569516
$ws.Range("J12").Value = 0.2221047250357578
$ws.Range("M12").Value = 0.9105153333333335
$ws.Range("N12").Value = 2.731546
$ws.Range("O12").Value = 0.03002598096739917
$ws.Range("P12").Value = 0.03049089185478593
$ws.Range("Q12").Value = 68.53983993510889
$ws.Range("R12").Value = 616.85855941598
$ws.Range("S12").Value = 0.00613770814551237
$ws.Range("T12").Value = 0.006772171151502254
$ws.Range("G13").Value = 75.27587666666666
$ws.Range("H13").Value = 225.82763
$ws.Range("I13").Value = 0.2044132430569516
$ws.Range("J13").Value = 0.2221047250357578
$ws.Range("O13").Value = 0.2937114969858886
$ws.Range("P13").Value = 0.2982592142727163
$ws.Range("Q13").Value = 670.4506677857177
$ws.Range("R13").Value = 6034.056010071459
$ws.Range("S13").Value = 0.06003851962199756
$ws.Range("T13").Value = 0.06624478077542283
$ws.Range("G14").Value = 75.27587666666666
$ws.Range("H14").Value = 225.82763
$ws.Range("I14").Value = 0.2044132430569516
$ws.Range("J14").Value = 0.2221047250357578
$ws.Range("M14").Value = 9.030046333333333
$ws.Range("N14").Value = 27.090139
$ws.Range("O14").Value = 0.2977830129963756
$ws.Range("P14").Value = 0.3023937720910131
$ws.Range("Q14").Value = 679.7446540822855
$ws.Range("R14").Value = 6117.70188674057
$ws.Range("S14").Value = 0.06087079141385952
$ws.Range("T14").Value = 0.06716308560280006
$ws.Range("G15").Value = 75.27587666666666
$ws.Range("H15").Value = 225.82763
$ws.Range("I15").Value = 0.2044132430569516
$ws.Range("J15").Value = 0.2221047250357578
$ws.Range("M15").Value = 1.38711
$ws.Range("N15").Value = 2.77422
$ws.Range("O15").Value = 0.04574259975086167
$ws.Range("P15").Value = 0.03096724053022875
$ws.Range("Q15").Value = 104.4159212831
$ws.Range("R15").Value = 626.4955276986
$ws.Range("S15").Value = 0.00935039316092974
$ws.Range("T15").Value = 0.006877970443082629
$ws.Range("G16").Value = 75.27587666666666
$ws.Range("H16").Value = 225.82763
$ws.Range("I16").Value = 0.2044132430569516
$ws.Range("J16").Value = 0.2221047250357578
$ws.Range("M16").Value = 10.089997
$ws.Range("N16").Value = 30.269991
$ws.Range("O16").Value = 0.3327369092994751
$ws.Range("P16").Value = 0.3378888812512559
$ws.Range("Q16").Value = 759.5333697390365
$ws.Range("R16").Value = 6835.800327651329
$ws.Range("S16").Value = 0.06801583071465248
$ws.Range("T16").Value = 0.07504671706295001
$ws.Range("G17").Value = 87.99833699999999
$ws.Range("H17").Value = 175.996674
$ws.Range("I17").Value = 0.2389613545046087
$ws.Range("J17").Value = 0.1730952624618072
$ws.Range("M17").Value = 0.9105153333333335
$ws.Range("N17").Value = 2.731546
$ws.Range("O17").Value = 0.03002598096739917
$ws.Range("P17").Value = 0.03049089185478593
$ws.Range("Q17").Value = 80.123835146334
$ws.Range("R17").Value = 480.743010878004
$ws.Range("S17").Value = 0.007175049082299307
$ws.Range("T17").Value = 0.00527782892829875
$ws.Range("G18").Value = 87.99833699999999
$ws.Range("H18").Value = 175.996674
$ws.Range("I18").Value = 0.2389613545046087
$ws.Range("J18").Value = 0.1730952624618072
$ws.Range("O18").Value = 0.2937114969858886
$ws.Range("P18").Value = 0.2982592142727163
$ws.Range("Q18").Value = 783.7642870230178
$ws.Range("R18").Value = 4702.585722138107
$ws.Range("S18").Value = 0.07018569715332423
$ws.Range("T18").Value = 0.05162725697618824
$ws.Range("G19").Value = 87.99833699999999
$ws.Range("H19").Value = 175.996674
$ws.Range("I19").Value = 0.2389613545046087
$ws.Range("J19").Value = 0.1730952624618072
$ws.Range("M19").Value = 9.030046333333333
$ws.Range("N19").Value = 27.090139
$ws.Range("O19").Value = 0.2977830129963756
$ws.Range("P19").Value = 0.3023937720910131
$ws.Range("Q19").Value = 794.6290603662809
$ws.Range("R19").Value = 4767.774362197686
$ws.Range("S19").Value = 0.07115863213407742
$ws.Range("T19").Value = 0.05234292934690983
$ws.Range("G20").Value = 87.99833699999999
$ws.Range("H20").Value = 175.996674
$ws.Range("I20").Value = 0.2389613545046087
$ws.Range("J20").Value = 0.1730952624618072
$ws.Range("M20").Value = 1.38711
$ws.Range("N20").Value = 2.77422
$ws.Range("O20").Value = 0.04574259975086167
$ws.Range("P20").Value = 0.03096724053022875
$ws.Range("Q20").Value = 122.06337323607
$ws.Range("R20").Value = 488.25349294428
$ws.Range("S20").Value = 0.01093071359502808
$ws.Range("T20").Value = 0.00536028262729786
$ws.Range("G21").Value = 87.99833699999999
$ws.Range("H21").Value = 175.996674
$ws.Range("I21").Value = 0.2389613545046087
$ws.Range("J21").Value = 0.1730952624618072
$ws.Range("M21").Value = 10.089997
$ws.Range("N21").Value = 30.269991
$ws.Range("O21").Value = 0.3327369092994751
$ws.Range("P21").Value = 0.3378888812512559
$ws.Range("Q21").Value = 887.9029563349887
$ws.Range("R21").Value = 5327.417738009933
$ws.Range("S21").Value = 0.07951126253987971
$ws.Range("T21").Value = 0.05848696458311256
$ws.Range("G22").Value = 88.99412
$ws.Range("H22").Value = 266.98236
$ws.Range("I22").Value = 0.24166542440621
$ws.Range("J22").Value = 0.2625809944389785
$ws.Range("M22").Value = 0.9105153333333335
$ws.Range("N22").Value = 2.731546
$ws.Range("O22").Value = 0.03002598096739917
$ws.Range("P22").Value = 0.03049089185478593
$ws.Range("Q22").Value = 81.03051083650668
$ws.Range("R22").Value = 729.27459752856
$ws.Range("S22").Value = 0.007256241433699304
$ws.Range("T22").Value = 0.008006328704561036
$ws.Range("G23").Value = 88.99412
$ws.Range("H23").Value = 266.98236
$ws.Range("I23").Value = 0.24166542440621
$ws.Range("J23").Value = 0.2625809944389785
$ws.Range("O23").Value = 0.2937114969858886
$ws.Range("P23").Value = 0.2982592142727163
$ws.Range("Q23").Value = 792.6333086390131
$ws.Range("R23").Value = 7133.699777751119
$ws.Range("S23").Value = 0.07097991357207803
$ws.Range("T23").Value = 0.07831720108431822
$ws.Range("G24").Value = 88.99412
$ws.Range("H24").Value = 266.98236
$ws.Range("I24").Value = 0.24166542440621
$ws.Range("J24").Value = 0.2625809944389785
$ws.Range("M24").Value = 9.030046333333333
$ws.Range("N24").Value = 27.090139
$ws.Range("O24").Value = 0.2977830129963756
$ws.Range("P24").Value = 0.3023937720910131
$ws.Range("Q24").Value = 803.6210269942266
$ws.Range("R24").Value = 7232.589242948039
$ws.Range("S24").Value = 0.07196385821672906
$ws.Range("T24").Value = 0.07940285738781204
$ws.Range("G25").Value = 88.99412
$ws.Range("H25").Value = 266.98236
$ws.Range("I25").Value = 0.24166542440621
$ws.Range("J25").Value = 0.2625809944389785
$ws.Range("M25").Value = 1.38711
$ws.Range("N25").Value = 2.77422
$ws.Range("O25").Value = 0.04574259975086167
$ws.Range("P25").Value = 0.03096724053022875
$ws.Range("Q25").Value = 123.4446337932
$ws.Range("R25").Value = 740.6678027592
$ws.Range("S25").Value = 0.01105440478223538
$ws.Range("T25").Value = 0.008131408813458503
$ws.Range("G26").Value = 88.99412
$ws.Range("H26").Value = 266.98236
$ws.Range("I26").Value = 0.24166542440621
$ws.Range("J26").Value = 0.2625809944389785
$ws.Range("M26").Value = 10.089997
$ws.Range("N26").Value = 30.269991
$ws.Range("O26").Value = 0.3327369092994751
$ws.Range("P26").Value = 0.3378888812512559
$ws.Range("Q26").Value = 897.9504038176398
$ws.Range("R26").Value = 8081.553634358758
$ws.Range("S26").Value = 0.08041100640146825
$ws.Range("T26").Value = 0.0887231984488287
